$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4997.091
$ws.Range("I38").Value = 133.8
$ws.Range("J38").Value = 9049.833000000001
$ws.Range("K38").Value = 401.4
$ws.Range("L38").Value = 27149.499
$ws.Range("M38").Value = -29.40000000000003
$ws.Range("N38").Value = -27893.499

$ws.Range("H43").Value = 5285.2856
$ws.Range("I43").Value = 4999.5
$ws.Range("J43").Value = 7000
$ws.Range("K43").Value = 4999.5
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = -4930.5
$ws.Range("N43").Value = -7138

$ws.Range("H48").Value = 2900.2856
$ws.Range("I48").Value = 1883.6666
$ws.Range("J48").Value = 9000
$ws.Range("K48").Value = 5650.9998
$ws.Range("L48").Value = 27000
$ws.Range("M48").Value = -5358.9998
$ws.Range("N48").Value = -27584

$ws.Range("H56").Value = 2900.2856
$ws.Range("I56").Value = 1883.6666
$ws.Range("J56").Value = 9000
$ws.Range("K56").Value = 5650.9998
$ws.Range("L56").Value = 27000
$ws.Range("M56").Value = -5116.9998
$ws.Range("N56").Value = -28068

$ws.Range("H58").Value = 991.8125
$ws.Range("I58").Value = 604.9286
$ws.Range("J58").Value = 3700
$ws.Range("K58").Value = 1814.7858
$ws.Range("L58").Value = 11100
$ws.Range("M58").Value = -1664.7858
$ws.Range("N58").Value = -11400

$ws.Range("H112").Value = 2514.5386
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2514.5386
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 7543.6158
$ws.Range("N112").Value = -9759.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2842.88
$ws.Range("I32").Value = 2487.4062
$ws.Range("J32").Value = 11374.25
$ws.Range("K32").Value = 2487.4062
$ws.Range("L32").Value = 11374.25
$ws.Range("M32").Value = -2200.4062
$ws.Range("N32").Value = -11948.25

$ws.Range("H45").Value = 20820.125
$ws.Range("I45").Value = 72921.664
$ws.Range("J45").Value = 3452.9443
$ws.Range("K45").Value = 72921.664
$ws.Range("L45").Value = 3452.9443
$ws.Range("M45").Value = -72544.664
$ws.Range("N45").Value = -4206.9443

$ws.Range("H61").Value = 5237.423
$ws.Range("I61").Value = 4287.3887
$ws.Range("J61").Value = 7375
$ws.Range("K61").Value = 4287.3887
$ws.Range("L61").Value = 7375
$ws.Range("M61").Value = -4075.3887
$ws.Range("N61").Value = -7799

$ws.Range("H136").Value = 5237.423
$ws.Range("I136").Value = 4287.3887
$ws.Range("J136").Value = 7375
$ws.Range("K136").Value = 12862.1661
$ws.Range("L136").Value = 22125
$ws.Range("M136").Value = -10312.1661
$ws.Range("N136").Value = -27225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2953.66
$ws.Range("I31").Value = 2757.0366
$ws.Range("J31").Value = 3849.389
$ws.Range("K31").Value = 2757.0366
$ws.Range("L31").Value = 3849.389
$ws.Range("M31").Value = -2462.0366
$ws.Range("N31").Value = -4439.389

$ws.Range("H34").Value = 2953.66
$ws.Range("I34").Value = 2757.0366
$ws.Range("J34").Value = 3849.389
$ws.Range("K34").Value = 2757.0366
$ws.Range("L34").Value = 3849.389
$ws.Range("M34").Value = -2555.0366
$ws.Range("N34").Value = -4253.389

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = ""

$ws.Range("H132").Value = 3368.1667
$ws.Range("I132").Value = 3264.5
$ws.Range("J132").Value = 4197.5
$ws.Range("K132").Value = 9793.5
$ws.Range("L132").Value = 12592.5
$ws.Range("M132").Value = -7263.5
$ws.Range("N132").Value = -17652.5

$ws.Range("H134").Value = 3752.853
$ws.Range("I134").Value = 3122.111
$ws.Range("J134").Value = 4462.4375
$ws.Range("K134").Value = 9366.332999999999
$ws.Range("L134").Value = 13387.3125
$ws.Range("M134").Value = -6831.332999999999
$ws.Range("N134").Value = -18457.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 154.75
$ws.Range("I23").Value = 41.25
$ws.Range("J23").Value = 268.25
$ws.Range("K23").Value = 123.75
$ws.Range("L23").Value = 804.75
$ws.Range("M23").Value = 111.25
$ws.Range("N23").Value = -1274.75

$ws.Range("H33").Value = 82.57692
$ws.Range("I33").Value = 75.85714
$ws.Range("J33").Value = 85.052635
$ws.Range("K33").Value = 455.14284
$ws.Range("L33").Value = 510.3158099999999
$ws.Range("M33").Value = -172.14284
$ws.Range("N33").Value = -1076.31581

$ws.Range("H113").Value = 2585.625
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2585.625
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7756.875
$ws.Range("N113").Value = -12096.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 28500
$ws.Range("I53").Value = 27000
$ws.Range("J53").Value = 30000
$ws.Range("K53").Value = 27000
$ws.Range("L53").Value = 30000
$ws.Range("M53").Value = -26369
$ws.Range("N53").Value = -31262

$ws.Range("H80").Value = 166669400
$ws.Range("I80").Value = 333335330
$ws.Range("J80").Value = 3498.3333
$ws.Range("K80").Value = 333335330
$ws.Range("L80").Value = 3498.3333
$ws.Range("M80").Value = -333334332
$ws.Range("N80").Value = -5494.3333

$ws.Range("H83").Value = 166669400
$ws.Range("I83").Value = 333335330
$ws.Range("J83").Value = 3498.3333
$ws.Range("K83").Value = 1666676650
$ws.Range("L83").Value = 17491.6665
$ws.Range("M83").Value = -1666671658
$ws.Range("N83").Value = -27475.6665

$ws.Range("H113").Value = 6296.391
$ws.Range("I113").Value = 6290.25
$ws.Range("J113").Value = 6337.3335
$ws.Range("K113").Value = 6290.25
$ws.Range("L113").Value = 6337.3335
$ws.Range("M113").Value = -4120.25
$ws.Range("N113").Value = -10677.3335

$ws.Range("H132").Value = 3757.8572
$ws.Range("I132").Value = 3391.762
$ws.Range("J132").Value = 4856.143
$ws.Range("K132").Value = 10175.286
$ws.Range("L132").Value = 14568.429
$ws.Range("M132").Value = -7645.286
$ws.Range("N132").Value = -19628.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1487.1351
$ws.Range("I46").Value = 2983.3333
$ws.Range("J46").Value = 1355.1177
$ws.Range("K46").Value = 2983.3333
$ws.Range("L46").Value = 1355.1177
$ws.Range("M46").Value = -2795.3333
$ws.Range("N46").Value = -1731.1177

$ws.Range("H61").Value = 3302.375
$ws.Range("I61").Value = 2131.5
$ws.Range("J61").Value = 4473.25
$ws.Range("K61").Value = 2131.5
$ws.Range("L61").Value = 4473.25
$ws.Range("M61").Value = -1929.5
$ws.Range("N61").Value = -4877.25

$ws.Range("H68").Value = 4002
$ws.Range("I68").Value = 3006.75
$ws.Range("J68").Value = 4997.25
$ws.Range("K68").Value = 3006.75
$ws.Range("L68").Value = 4997.25
$ws.Range("M68").Value = -2257.75
$ws.Range("N68").Value = -6495.25

$ws.Range("H71").Value = 4002
$ws.Range("I71").Value = 3006.75
$ws.Range("J71").Value = 4997.25
$ws.Range("K71").Value = 15033.75
$ws.Range("L71").Value = 24986.25
$ws.Range("M71").Value = -11289.75
$ws.Range("N71").Value = -32474.25

$ws.Range("H82").Value = 6135.5
$ws.Range("I82").Value = 5216
$ws.Range("J82").Value = 8894
$ws.Range("K82").Value = 5216
$ws.Range("L82").Value = 8894
$ws.Range("M82").Value = -4855
$ws.Range("N82").Value = -9616

$ws.Range("H85").Value = 6135.5
$ws.Range("I85").Value = 5216
$ws.Range("J85").Value = 8894
$ws.Range("K85").Value = 5216
$ws.Range("L85").Value = 8894
$ws.Range("M85").Value = -3968
$ws.Range("N85").Value = -11390

$ws.Range("H93").Value = 2253.125
$ws.Range("I93").Value = 2117
$ws.Range("J93").Value = 2661.5
$ws.Range("K93").Value = 2117
$ws.Range("L93").Value = 2661.5
$ws.Range("M93").Value = -869
$ws.Range("N93").Value = -5157.5

$ws.Range("H113").Value = 3302.375
$ws.Range("I113").Value = 2131.5
$ws.Range("J113").Value = 4473.25
$ws.Range("K113").Value = 2131.5
$ws.Range("L113").Value = 4473.25
$ws.Range("M113").Value = 38.5
$ws.Range("N113").Value = -8813.25

$ws.Range("H122").Value = 8801
$ws.Range("I122").Value = 9333.333000000001
$ws.Range("J122").Value = 8002.5
$ws.Range("K122").Value = 27999.999
$ws.Range("L122").Value = 24007.5
$ws.Range("M122").Value = -25549.999
$ws.Range("N122").Value = -28907.5

$ws.Range("H132").Value = 4816.2666
$ws.Range("I132").Value = 4156.125
$ws.Range("J132").Value = 5570.7144
$ws.Range("K132").Value = 12468.375
$ws.Range("L132").Value = 16712.1432
$ws.Range("M132").Value = -9938.375
$ws.Range("N132").Value = -21772.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11500
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -12748

$ws.Range("H65").Value = 11500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 57500
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -63740

$ws.Range("H81").Value = 6957
$ws.Range("I81").Value = 7599.75
$ws.Range("J81").Value = 6699.9
$ws.Range("K81").Value = 15199.5
$ws.Range("L81").Value = 13399.8
$ws.Range("M81").Value = -14138.5
$ws.Range("N81").Value = -15521.8

$ws.Range("H84").Value = 6957
$ws.Range("I84").Value = 7599.75
$ws.Range("J84").Value = 6699.9
$ws.Range("K84").Value = 75997.5
$ws.Range("L84").Value = 66999
$ws.Range("M84").Value = -70693.5
$ws.Range("N84").Value = -77607

$ws.Range("H96").Value = 5595.8
$ws.Range("I96").Value = 4567.2856
$ws.Range("J96").Value = 7995.6665
$ws.Range("K96").Value = 4567.2856
$ws.Range("L96").Value = 7995.6665
$ws.Range("M96").Value = -3194.2856
$ws.Range("N96").Value = -10741.6665

$ws.Range("H122").Value = 8335652
$ws.Range("I122").Value = 1852.2174
$ws.Range("J122").Value = 35718136
$ws.Range("K122").Value = 5556.6522
$ws.Range("L122").Value = 107154408
$ws.Range("M122").Value = -3106.6522
$ws.Range("N122").Value = -107159308

$ws.Range("H132").Value = 3808.2903
$ws.Range("I132").Value = 3534.2454
$ws.Range("J132").Value = 5422.1113
$ws.Range("K132").Value = 10602.7362
$ws.Range("L132").Value = 16266.3339
$ws.Range("M132").Value = -8072.736199999999
$ws.Range("N132").Value = -21326.3339

$ws.Range("H136").Value = 22728744
$ws.Range("I136").Value = 27028288
$ws.Range("J136").Value = 2582.5715
$ws.Range("K136").Value = 81084864
$ws.Range("L136").Value = 7747.7145
$ws.Range("M136").Value = -81082314
$ws.Range("N136").Value = -12847.7145
